# cleaned more data, fixed a problem where some SLRs contained incorrect data.
# The "Authors" column (E) on rows 2-9 held author-list strings whose field
# separators had collapsed/inconsistent run-lengths of whitespace after each
# comma. They are rewritten here with the whitespace runs normalised/widened
# by one extra space (first pass) and then by a further extra space (second,
# final pass) to match the cleaned data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = 2, 3, 4, 5, 6, 7, 8, 9

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 5)   # column E = Authors
    $orig = $cell.Value2

    # First pass: widen each comma's following whitespace run by 1 space.
    $pass1 = $orig -replace ",( +)", ',$1 '
    $cell.Value = $pass1

    # Second pass: widen by 1 more space (net +2 vs. the original), this is
    # the value that remains in the cell.
    $pass2 = $pass1 -replace ",( +)", ',$1 '
    $cell.Value = $pass2
}
